# edit.ps1 -- reproduces the "updated system overview image for readme" commit
#
# Changes applied:
#   1. Every cached "datetimeFigureOut" field (slide master, all 11 slide
#      layouts, notes master) is bumped from 10/02/2018 -> 13/02/2018.
#   2. Two shapes on slide 1 ("Rounded Rectangle 71" / "Storage Interface"
#      and the rotated "Left-Right Arrow 83" connecting it) are nudged /
#      resized slightly.
#   3. The "Pawsey Supercomputing Centre" caption textbox is retitled to
#      "GPU Cluster".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Date placeholder fields: 10/02/2018 -> 13/02/2018
#    (slide master + every custom layout + the notes master)
# ---------------------------------------------------------------------
function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
                if ($sh.TextFrame.TextRange.Text -eq "10/02/2018") {
                    $sh.TextFrame.TextRange.Text = "13/02/2018"
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholders $layout.Shapes
}

Update-DatePlaceholders $p.NotesMaster.Shapes

# ---------------------------------------------------------------------
# 2. Reposition / resize the two connector shapes on slide 1
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)

$storageInterface = $slide.Shapes.Item("Rounded Rectangle 71")
$storageInterface.Left = 432.6803149606299
$storageInterface.Top = 205.65969093937008

$connectorArrow = $slide.Shapes.Item("Left-Right Arrow 83")
$connectorArrow.Left = 455.8596954393701
$connectorArrow.Width = 37.03858377716535

# ---------------------------------------------------------------------
# 3. Retitle the caption textbox
# ---------------------------------------------------------------------
$caption = $slide.Shapes.Item("TextBox 59")
$caption.TextFrame.TextRange.Text = "GPU Cluster"
